# Apply updated cryptocurrency price/volume figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.413.90"
$ws.Range("E2").Value = "  +0.68%  "

$ws.Range("D3").Value = "3.525.65"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("D5").Value = "'597.41"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").Value = "'173.83"
$ws.Range("E6").Value = "  +2.64%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.595"
$ws.Range("E8").Value = "  +3.36%  "

$ws.Range("E9").Value = "  +7.98%  "

$ws.Range("E10").Value = "  +0.27%  "

$ws.Range("D11").Value = "'0.438"
$ws.Range("E11").Value = "  -0.33%  "

$ws.Range("D12").Value = "4.133.37"
$ws.Range("E12").Value = "  +0.59%  "

$ws.Range("E13").Value = "  +0.11%  "

$ws.Range("D14").Value = "'28.87"
$ws.Range("E14").Value = "  +2.22%  "

$ws.Range("E15").Value = "  +1.79%  "

$ws.Range("D16").Value = "67.342.43"
$ws.Range("E16").Value = "  +0.66%  "

$ws.Range("D17").Value = "3.519.08"
$ws.Range("E17").Value = "  +0.23%  "

$ws.Range("D18").Value = "'6.35"
$ws.Range("E18").Value = "  +0.39%  "

$ws.Range("D19").Value = "'14.31"
$ws.Range("E19").Value = "  +1.78%  "

$ws.Range("D20").Value = "'397.84"
$ws.Range("E20").Value = "  +0.69%  "

$ws.Range("D21").Value = "'8.00"
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("D22").Value = "'73.47"
$ws.Range("E22").Value = "  +0.37%  "

$ws.Range("D23").Value = "'0.541"
$ws.Range("E23").Value = "  +1.76%  "

$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("D25").Value = "'0.0000124"
$ws.Range("E25").Value = "  -3.41%  "

$ws.Range("D26").Value = "'10.30"
$ws.Range("E26").Value = "  +2.19%  "

$ws.Range("E27").Value = "  -0.62%  "

$ws.Range("E28").Value = "  -0.27%  "

$ws.Range("D29").Value = "'6.30"
$ws.Range("E29").Value = "  -1.16%  "

$ws.Range("E30").Value = "  -0.36%  "

$ws.Range("E31").Value = "  +0.53%  "

$ws.Range("D32").Value = "'24.18"
$ws.Range("E32").Value = "  +2.55%  "

$ws.Range("D33").Value = "'7.41"
$ws.Range("E33").Value = "  -0.72%  "

$ws.Range("E34").Value = "  +2.76%  "

$ws.Range("D35").Value = "'163.33"
$ws.Range("E35").Value = "  +0.46%  "

$ws.Range("D36").Value = "'0.897"
$ws.Range("E36").Value = "  -0.52%  "

$ws.Range("E37").Value = "  -1.16%  "

$ws.Range("D38").Value = "'6.96"
$ws.Range("E38").Value = "  +4.16%  "

$ws.Range("D39").Value = "'4.73"
$ws.Range("E39").Value = "  +1.19%  "

$ws.Range("D40").Value = "'27.74"
$ws.Range("E40").Value = "  +4.35%  "

$ws.Range("D41").Value = "'0.0747"
$ws.Range("E41").Value = "  -1.03%  "

$ws.Range("D42").Value = "'26.48"
$ws.Range("E42").Value = "  +0.46%  "

$ws.Range("D43").Value = "'2.63"
$ws.Range("E43").Value = "  +2.97%  "

$ws.Range("D44").Value = "2.804.88"
$ws.Range("E44").Value = "  -1.28%  "

$ws.Range("D45").Value = "'42.93"
$ws.Range("E45").Value = "  -1.32%  "

$ws.Range("E46").Value = "  -2.46%  "

$ws.Range("D47").Value = "'341.23"
$ws.Range("E47").Value = "  -2.17%  "

$ws.Range("E48").Value = "  +1.69%  "

$ws.Range("D49").Value = "'33.88"
$ws.Range("E49").Value = "  +0.27%  "

$ws.Range("D50").Value = "'6.55"
$ws.Range("E50").Value = "  -0.04%  "

$ws.Range("D51").Value = "'0.854"
$ws.Range("E51").Value = "  -0.75%  "
